$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eigen_Terrasse")

# Insert a new row 9 for the "Dachrinne" option (shifts old rows 9-12 down to 10-13)
$ws.Rows.Item(9).Insert()

# Insert another row at 13 to leave a gap before the final "Preis" row (old row 12 -> now 14)
$ws.Rows.Item(13).Insert()

# Fill the new row 9 with the Dachrinne (gutter) option fields
$ws.Cells.Item(9, 1).Value = "Auswahl"
$ws.Cells.Item(9, 2).Value = "Dachrinne"
$ws.Cells.Item(9, 3).Value = "P_Rinne"
$ws.Cells.Item(9, 4).Value = "Nein:0, Ja:205"

# Update the final price formula (row 14, column E) to add the gutter price term
$ws.Cells.Item(14, 5).Value = "( (L * P_Trager * P_OF) + (N_Col * H * 90* P_OF) + (N_Spar * B * 110* P_OF) + (L * B * P_Dach) + (L * P_Wand) + (L * P_Rinne) * (1 - (Rabatt / 100))"

# Update selection to match the new active cell
$ws.Range("E14").Select()
